$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XLNames")

# Delete row 2 (CurrencyMasterData.xlsx), shifting the remaining rows up
$ws.Rows.Item(2).Delete()

# Update the active selection on the sheet
$ws.Range("G14").Select()
